$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 47.3758915439561
$ws.Range("K2").Value = 46.7829418487092
$ws.Range("N2").Value = 46.7614395656938

$ws.Range("B3").Value = 17.727547636372
$ws.Range("K3").Value = 17.7374744009919
$ws.Range("N3").Value = 17.8066357625234

$ws.Range("B4").Value = 4.66246711093616
$ws.Range("K4").Value = 4.54122555886245
$ws.Range("M4").Value = 4.98006182286059
$ws.Range("N4").Value = 4.50479434655862

$ws.Range("B5").Value = 2904.37011684753
$ws.Range("K5").Value = 2791.49348527449
$ws.Range("M5").Value = 3291.56572730894
$ws.Range("N5").Value = 2790.15624910416

$ws.Range("B6").Value = 5.11353961710558
$ws.Range("K6").Value = 4.65867992158079
$ws.Range("N6").Value = 4.90224996427908

$ws.Range("B7").Value = 0.557752965581924
$ws.Range("D7").Value = 0.528351073871356
$ws.Range("K7").Value = 0.396376822894605
$ws.Range("N7").Value = 0.583518233303355

$ws.Range("B8").Value = 0.429454246552531
$ws.Range("K8").Value = 0.592127159062159
$ws.Range("N8").Value = 0.408149050854568

$ws.Range("B9").Value = 0.351197281372668
$ws.Range("K9").Value = 0.475466970270242
$ws.Range("N9").Value = 0.406100171031648

$ws.Range("B10").Value = 0.398173738402363
$ws.Range("K10").Value = 0.279970941185883
$ws.Range("N10").Value = 0.426839676660333
